# Sprint2 Backlog/Burndown update — "Added Coverage for EmailService"
#
# The Demo-notes follow-up items (rows 26-40) get their estimate/remaining
# numbers filled in and several of them get re-assigned from Ahmad/Nick to
# Trinidad. Row 47's "Amount Remaining After Week 2" also gets filled in
# with 0. The Estimate Totals row (48) and the burndown chart both recalc
# automatically off of the SUM formulas already in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26 - (minor) error message centered on login...
$ws.Range("D26").Value = 0
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = "Trinidad"

# Row 27 - manager should be able to add/remove employee from their own group
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 0
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = "Trinidad"

# Row 28 - user should not be set to a system wide manager role...
$ws.Range("C28").Value = 1
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = "Trinidad"

# Row 29 - should not allow deleting a user that is set as manager...
$ws.Range("C29").Value = 2
$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = "Trinidad"

# Row 32 - ensure forms indicate required fields and other expectations
$ws.Range("C32").Value = 2
$ws.Range("D32").Value = 0
$ws.Range("E32").Value = 0
$ws.Range("F32").Value = "Trinidad"

# Row 33 - most error messages are consistent with a few exceptions...
$ws.Range("C33").Value = 2
$ws.Range("D33").Value = 0
$ws.Range("E33").Value = 0
$ws.Range("F33").Value = "Trinidad"

# Row 34 - add filtering/paging/scrolling to support working with large lists
$ws.Range("C34").Value = 2
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("F34").Value = "Trinidad"

# Row 36 - add confirmation check when deleting group
$ws.Range("E36").Value = 0
$ws.Range("F36").Value = "Trinidad"

# Row 38 - employee should be able to view all group/project
$ws.Range("C38").Value = 2
$ws.Range("D38").Value = 0
$ws.Range("E38").Value = 0
$ws.Range("F38").Value = "Trinidad"

# Row 39 - manager should be able to edit/delete their own group/project
$ws.Range("C39").Value = 2
$ws.Range("D39").Value = 0
$ws.Range("E39").Value = 0
$ws.Range("F39").Value = "Trinidad"

# Row 47 - Sprint Backlog: Amount Remaining After Week 2 filled in
$ws.Range("E47").Value = 0

# Move the active selection (matches the saved view in the edited workbook)
$ws.Range("B7").Select()
